# Commit: Update countries & provincias Spain
# Applies the sharedStrings reorder + refreshed COVID stats captured in the diff,
# by writing the resulting cell values directly (country names + numeric columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 21:29"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 2924771
$ws.Range("C4").Value = 34183
$ws.Range("D4").Value = 1251542
$ws.Range("E4").Value = 1541006
$ws.Range("G4").Value = 159
$ws.Range("H4").Value = 132223

# Row 7: India -> India
$ws.Range("B7").Value = 672702
$ws.Range("C7").Value = 22813
$ws.Range("D7").Value = 408645
$ws.Range("E7").Value = 244778

# Row 10: Chile -> Chile
$ws.Range("B10").Value = 291847
$ws.Range("C10").Value = 3758
$ws.Range("D10").Value = 257445
$ws.Range("E10").Value = 28210
$ws.Range("G10").Value = 141
$ws.Range("H10").Value = 6192

# Row 49: Rumania -> Israel
$ws.Range("A49").Value = "Israel"
$ws.Range("B49").Value = 29032
$ws.Range("C49").Value = 977
$ws.Range("D49").Value = 17773
$ws.Range("E49").Value = 10929
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 330

# Row 50: Barein -> Rumania
$ws.Range("A50").Value = "Rumania"
$ws.Range("B50").Value = 28582
$ws.Range("C50").Value = 416
$ws.Range("D50").Value = 19854
$ws.Range("E50").Value = 6997
$ws.Range("G50").Value = 23
$ws.Range("H50").Value = 1731

# Row 51: Israel -> Barein
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 28410
$ws.Range("D51").Value = 23318
$ws.Range("E51").Value = 4996
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 96

# Row 71: Sudan -> Uzbekistan
$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("B71").Value = 9708
$ws.Range("C71").Value = 312
$ws.Range("D71").Value = 6425
$ws.Range("E71").Value = 3252
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 31

# Row 72: Uzbekistan -> Sudan
$ws.Range("A72").Value = "Sudan"
$ws.Range("B72").Value = 9663
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 4624
$ws.Range("E72").Value = 4435
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 604

# Row 94: Luxemburgo -> Costa Rica
$ws.Range("A94").Value = "Costa Rica"
$ws.Range("B94").Value = 4621
$ws.Range("C94").Value = 310
$ws.Range("D94").Value = 1721
$ws.Range("E94").Value = 2882
$ws.Range("H94").Value = 18

# Row 95: Costa Rica -> Luxemburgo
$ws.Range("A95").Value = "Luxemburgo"
$ws.Range("B95").Value = 4476
$ws.Range("C95").Value = 29
$ws.Range("D95").Value = 4016
$ws.Range("E95").Value = 350
$ws.Range("H95").Value = 110

# Row 98: Estado de Palestina -> Estado de Palestina
$ws.Range("B98").Value = 3835
$ws.Range("C98").Value = 501
$ws.Range("E98").Value = 3360

# Row 111: Mali -> Mali
$ws.Range("B111").Value = 2303
$ws.Range("C111").Value = 18
$ws.Range("D111").Value = 1516
$ws.Range("E111").Value = 669
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 118

# Row 131: Jordania -> Jordania
$ws.Range("D131").Value = 902
$ws.Range("E131").Value = 238

# Row 136: Burkina Faso -> Burkina Faso
$ws.Range("D136").Value = 854
$ws.Range("E136").Value = 73

# Row 143: Republica del Chad -> Republica del Chad
$ws.Range("D143").Value = 787
$ws.Range("E143").Value = 10

# Row 209: Groenlandia -> Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"

# Row 210: Islas Malvinas -> Groenlandia
$ws.Range("A210").Value = "Groenlandia"
